$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Insert two new columns at E (zve_abg_nokfb, zve_abg_kfb).
#    This pushes old E (gross_e5) -> G, old F (year) -> H,
#    old H (tax_nokfb) -> J, old K/L/M/N (G/M/S/R consts) -> M/N/O/P,
#    old S (helper col) -> U
# -----------------------------------------------------------------
$ws.Range("E1:F1").EntireColumn.Insert()

# -----------------------------------------------------------------
# 2) Insert three more new columns at M (after shift, this is the
#    spot right before the old "G" boundary-constant column), to make
#    room for tax_kfb, tax_abg_nokfb, tax_abg_kfb, abgst, abgst_tu
#    (two of the five new header columns, K & L, are already free
#    blank spacer columns left over from the old layout).
#    This pushes the boundary-constant block (currently M/N/O/P) to
#    P/Q/R/S, and the old helper column (currently U) to X.
# -----------------------------------------------------------------
$ws.Range("M1:O1").EntireColumn.Insert()

# -----------------------------------------------------------------
# Header row (row 1) - written in the same order the original author
# typed them, so the shared-string table comes out in the same order.
# -----------------------------------------------------------------
$ws.Range("K1").Value = "tax_kfb"
$ws.Range("N1").Value = "abgst"
$ws.Range("O1").Value = "abgst_tu"
$ws.Range("M1").Value = "tax_abg_kfb"
$ws.Range("L1").Value = "tax_abg_nokfb"
$ws.Range("E1").Value = "zve_abg_nokfb"
$ws.Range("F1").Value = "zve_abg_kfb"

# -----------------------------------------------------------------
# New value columns E (zve_abg_nokfb) & F (zve_abg_kfb), rows 2-6
# -----------------------------------------------------------------
$ws.Range("E2").Value = 5500
$ws.Range("F2").Value = 5500
$ws.Range("E3").Value = 12000
$ws.Range("F3").Value = 12000
$ws.Range("E4").Value = 21000
$ws.Range("F4").Value = 21000
$ws.Range("E5").Value = 52500
$ws.Range("F5").Value = 44500
$ws.Range("E6").Value = 200000
$ws.Range("F6").Value = 180000

# -----------------------------------------------------------------
# New column G (gross_e5 difference amounts), rows 2-6
# -----------------------------------------------------------------
$ws.Range("G2").Value = 500
$ws.Range("G3").Value = 2000
$ws.Range("G4").Value = 1000
$ws.Range("G5").Value = 2500
$ws.Range("G6").Value = 0

# -----------------------------------------------------------------
# Column D (zve_kfb): rows 2, 5, 6 lose their formula and become
# plain values; row 3's shared formula range shrinks to D3:D4.
# -----------------------------------------------------------------
$ws.Range("D2").Value = 5000
$ws.Range("D3").Formula = "=C3"
$ws.Range("D4").Formula = "=C4"
$ws.Range("D5").Value = 42000
$ws.Range("D6").Value = 180000

# -----------------------------------------------------------------
# Column J (tax_nokfb, old H) formulas - rows 3 & 4 get corrected
# constants; all rows now reference the shifted boundary columns
# P/Q/R/S instead of the old K/L/M/N.
# -----------------------------------------------------------------
$ws.Range("J2").Formula = "=TRUNC((C2>9000)*(C2<13996)*(997.8*(C2-9000)/10000+1400)*(C2-9000)/10000+(C2>13996)*(C2<54949)*((220.13*(C2-13996)/10000+2397)*(C2-13996)/10000+948.49)+(C2>54950)*(C2<260532)*(0.42*C2-8621.75)+(C2>260532)*(0.45*C2-16437.7))"
$ws.Range("J3").Formula = "=TRUNC((C3>P3)*(C3<Q3)*(974.58*(C3-P3)/10000+1400)*(C3-P3)/10000+(C3>Q3)*(C3<R3)*((228.74*(C3-Q3)/10000+2397)*(C3-Q3)/10000+971)+(C3>R3)*(C3<S3)*(0.42*C3-8239)+(C3>S3)*(0.45*C3-15761))"
$ws.Range("J4").Formula = "=TRUNC((C4>P4)*(C4<Q4)*(939.68*(C4-P4)/10000+1400)*(C4-P4)/10000+(C4>Q4)*(C4<R4)*((228.74*(C4-Q4)/10000+2397)*(C4-Q4)/10000+1007)+(C4>R4)*(C4<S4)*(0.42*C4-8064)+(C4>S4)*(0.45*C4-15576))"
$ws.Range("J5").Formula = "=TRUNC((C5>P5)*(C5<Q5)*(883.74*(C5-P5)/10000+1500)*(C5-P5)/10000+(C5>Q5)*(C5<R5)*((228.74*(C5-Q5)/10000+2397)*(C5-Q5)/10000+989)+(C5>R5)*(C5<S5)*(0.42*C5-7914)+(C5>S5)*(0.45*C5-15414))"
$ws.Range("J6").Formula = "=TRUNC((C6>9000)*(C6<13996)*(997.8*(C6-9000)/10000+1400)*(C6-9000)/10000+(C6>13996)*(C6<54949)*((220.13*(C6-13996)/10000+2397)*(C6-13996)/10000+948.49)+(C6>54950)*(C6<260532)*(0.42*C6-8621.75)+(C6>260532)*(0.45*C6-16437.7))"

# -----------------------------------------------------------------
# Column K (tax_kfb) - same TRUNC pattern, applied to column D
# -----------------------------------------------------------------
$ws.Range("K2").Formula = "=TRUNC((D2>9000)*(D2<13996)*(997.8*(D2-9000)/10000+1400)*(D2-9000)/10000+(D2>13996)*(D2<54949)*((220.13*(D2-13996)/10000+2397)*(D2-13996)/10000+948.49)+(D2>54950)*(D2<260532)*(0.42*D2-8621.75)+(D2>260532)*(0.45*D2-16437.7))"
$ws.Range("K3").Formula = "=TRUNC((D3>P3)*(D3<Q3)*(974.58*(D3-P3)/10000+1400)*(D3-P3)/10000+(D3>Q3)*(D3<R3)*((228.74*(D3-Q3)/10000+2397)*(D3-Q3)/10000+971)+(D3>R3)*(D3<S3)*(0.42*D3-8239)+(D3>S3)*(0.45*D3-15761))"
$ws.Range("K4").Formula = "=TRUNC((D4>P4)*(D4<Q4)*(939.68*(D4-P4)/10000+1400)*(D4-P4)/10000+(D4>Q4)*(D4<R4)*((228.74*(D4-Q4)/10000+2397)*(D4-Q4)/10000+1007)+(D4>R4)*(D4<S4)*(0.42*D4-8064)+(D4>S4)*(0.45*D4-15576))"
$ws.Range("K5").Formula = "=TRUNC((D5>P5)*(D5<Q5)*(883.74*(D5-P5)/10000+1500)*(D5-P5)/10000+(D5>Q5)*(D5<R5)*((228.74*(D5-Q5)/10000+2397)*(D5-Q5)/10000+989)+(D5>R5)*(D5<S5)*(0.42*D5-7914)+(D5>S5)*(0.45*D5-15414))"
$ws.Range("K6").Formula = "=TRUNC((D6>9000)*(D6<13996)*(997.8*(D6-9000)/10000+1400)*(D6-9000)/10000+(D6>13996)*(D6<54949)*((220.13*(D6-13996)/10000+2397)*(D6-13996)/10000+948.49)+(D6>54950)*(D6<260532)*(0.42*D6-8621.75)+(D6>260532)*(0.45*D6-16437.7))"

# -----------------------------------------------------------------
# Column L (tax_abg_nokfb) - same TRUNC pattern, applied to column E
# -----------------------------------------------------------------
$ws.Range("L2").Formula = "=TRUNC((E2>9000)*(E2<13996)*(997.8*(E2-9000)/10000+1400)*(E2-9000)/10000+(E2>13996)*(E2<54949)*((220.13*(E2-13996)/10000+2397)*(E2-13996)/10000+948.49)+(E2>54950)*(E2<260532)*(0.42*E2-8621.75)+(E2>260532)*(0.45*E2-16437.7))"
$ws.Range("L3").Formula = "=TRUNC((E3>P3)*(E3<Q3)*(974.58*(E3-P3)/10000+1400)*(E3-P3)/10000+(E3>Q3)*(E3<R3)*((228.74*(E3-Q3)/10000+2397)*(E3-Q3)/10000+971)+(E3>R3)*(E3<S3)*(0.42*E3-8239)+(E3>S3)*(0.45*E3-15761))"
$ws.Range("L4").Formula = "=TRUNC((E4>P4)*(E4<Q4)*(939.68*(E4-P4)/10000+1400)*(E4-P4)/10000+(E4>Q4)*(E4<R4)*((228.74*(E4-Q4)/10000+2397)*(E4-Q4)/10000+1007)+(E4>R4)*(E4<S4)*(0.42*E4-8064)+(E4>S4)*(0.45*E4-15576))"
$ws.Range("L5").Formula = "=TRUNC((E5>P5)*(E5<Q5)*(883.74*(E5-P5)/10000+1500)*(E5-P5)/10000+(E5>Q5)*(E5<R5)*((228.74*(E5-Q5)/10000+2397)*(E5-Q5)/10000+989)+(E5>R5)*(E5<S5)*(0.42*E5-7914)+(E5>S5)*(0.45*E5-15414))"
$ws.Range("L6").Formula = "=TRUNC((E6>9000)*(E6<13996)*(997.8*(E6-9000)/10000+1400)*(E6-9000)/10000+(E6>13996)*(E6<54949)*((220.13*(E6-13996)/10000+2397)*(E6-13996)/10000+948.49)+(E6>54950)*(E6<260532)*(0.42*E6-8621.75)+(E6>260532)*(0.45*E6-16437.7))"

# -----------------------------------------------------------------
# Column M (tax_abg_kfb) - same TRUNC pattern, applied to column F
# (row 3's formula reproduces the original author's typo exactly)
# -----------------------------------------------------------------
$ws.Range("M2").Formula = "=TRUNC((F2>9000)*(F2<13996)*(997.8*(F2-9000)/10000+1400)*(F2-9000)/10000+(F2>13996)*(F2<54949)*((220.13*(F2-13996)/10000+2397)*(F2-13996)/10000+948.49)+(F2>54950)*(F2<260532)*(0.42*F2-8621.75)+(F2>260532)*(0.45*F2-16437.7))"
$ws.Range("M3").Formula = "=TRUNC((F3>P3)*(F3<Q3)*(974.58*(F3-P3)/10000+1400)*(F3-P3)/10000+(F3>Q3)*(F3<S3)*((228.74*(F3-Q3)/10000+2397)*(F3-R3)/10000+971)+(F3>S3)*(F3<S3)*(0.42*F3-8239)+(F3>S3)*(0.45*F3-15761))"
$ws.Range("M4").Formula = "=TRUNC((F4>P4)*(F4<Q4)*(939.68*(F4-P4)/10000+1400)*(F4-P4)/10000+(F4>Q4)*(F4<R4)*((228.74*(F4-Q4)/10000+2397)*(F4-Q4)/10000+1007)+(F4>R4)*(F4<S4)*(0.42*F4-8064)+(F4>S4)*(0.45*F4-15576))"
$ws.Range("M5").Formula = "=TRUNC((F5>P5)*(F5<Q5)*(883.74*(F5-P5)/10000+1500)*(F5-P5)/10000+(F5>Q5)*(F5<R5)*((228.74*(F5-Q5)/10000+2397)*(F5-Q5)/10000+989)+(F5>R5)*(F5<S5)*(0.42*F5-7914)+(F5>S5)*(0.45*F5-15414))"
$ws.Range("M6").Formula = "=TRUNC((F6>9000)*(F6<13996)*(997.8*(F6-9000)/10000+1400)*(F6-9000)/10000+(F6>13996)*(F6<54949)*((220.13*(F6-13996)/10000+2397)*(F6-13996)/10000+948.49)+(F6>54950)*(F6<260532)*(0.42*F6-8621.75)+(F6>260532)*(0.45*F6-16437.7))"

# -----------------------------------------------------------------
# Column N (abgst) and O (abgst_tu), rows 2-6
# Rows 3-6 share two formula groups (si=2 for N, si=3 for O).
# -----------------------------------------------------------------
$ws.Range("N2").Formula = "=MAX((G2-801)*0.25,0)"
$ws.Range("O2").Formula = "=N2"
$ws.Range("N3:N6").Formula = "=MAX((G3-801)*0.25,0)"
$ws.Range("O3:O6").Formula = "=N3"

# -----------------------------------------------------------------
# SheetView settings
# -----------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("P7").Select()
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("H1").Column
